$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.417859666666667
$ws.Range("H2").Value = 4.253579
$ws.Range("I2").Value = 0.1472651073415806
$ws.Range("J2").Value = 0.1472651073415806
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.244885
$ws.Range("N2").Value = 72.734655
$ws.Range("O2").Value = 0.6895205882382217
$ws.Range("P2").Value = 0.6895205882382218
$ws.Range("Q2").Value = 34.37584456447167
$ws.Range("R2").Value = 309.382601080245
$ws.Range("S2").Value = 0.1015423234411315
$ws.Range("T2").Value = 0.1015423234411316

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.417859666666667
$ws.Range("H3").Value = 4.253579
$ws.Range("I3").Value = 0.1472651073415806
$ws.Range("J3").Value = 0.1472651073415806
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7651789999999999
$ws.Range("N3").Value = 2.295537
$ws.Range("O3").Value = 0.02176156637523891
$ws.Range("P3").Value = 0.02176156637523891
$ws.Range("Q3").Value = 1.084916441880333
$ws.Range("R3").Value = 9.764247976923
$ws.Range("S3").Value = 0.003204719408170489
$ws.Range("T3").Value = 0.00320471940817049

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.417859666666667
$ws.Range("H4").Value = 4.253579
$ws.Range("I4").Value = 0.1472651073415806
$ws.Range("J4").Value = 0.1472651073415806
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6351283333333334
$ws.Range("N4").Value = 1.905385
$ws.Range("O4").Value = 0.01806294655581008
$ws.Range("P4").Value = 0.01806294655581008
$ws.Range("Q4").Value = 0.9005228469905557
$ws.Range("R4").Value = 8.104705622915001
$ws.Range("S4").Value = 0.002660041763446605
$ws.Range("T4").Value = 0.002660041763446606

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.417859666666667
$ws.Range("H5").Value = 4.253579
$ws.Range("I5").Value = 0.1472651073415806
$ws.Range("J5").Value = 0.1472651073415806
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.516752666666667
$ws.Range("N5").Value = 28.550258
$ws.Range("O5").Value = 0.2706548988307292
$ws.Range("P5").Value = 0.2706548988307293
$ws.Range("Q5").Value = 13.49341976370911
$ws.Range("R5").Value = 121.440777873382
$ws.Range("S5").Value = 0.03985802272883198
$ws.Range("T5").Value = 0.03985802272883199

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.189892666666666
$ws.Range("H6").Value = 18.569678
$ws.Range("I6").Value = 0.6429093297593833
$ws.Range("J6").Value = 0.6429093297593833
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.244885
$ws.Range("N6").Value = 72.734655
$ws.Range("O6").Value = 0.6895205882382217
$ws.Range("P6").Value = 0.6895205882382218
$ws.Range("Q6").Value = 150.0732358656767
$ws.Range("R6").Value = 1350.65912279109
$ws.Range("S6").Value = 0.4432992192395309
$ws.Range("T6").Value = 0.4432992192395309

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.189892666666666
$ws.Range("H7").Value = 18.569678
$ws.Range("I7").Value = 0.6429093297593833
$ws.Range("J7").Value = 0.6429093297593833
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7651789999999999
$ws.Range("N7").Value = 2.295537
$ws.Range("O7").Value = 0.02176156637523891
$ws.Range("P7").Value = 0.02176156637523891
$ws.Range("Q7").Value = 4.736375880787333
$ws.Range("R7").Value = 42.627382927086
$ws.Range("S7").Value = 0.01399071405281918
$ws.Range("T7").Value = 0.01399071405281918

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.189892666666666
$ws.Range("H8").Value = 18.569678
$ws.Range("I8").Value = 0.6429093297593833
$ws.Range("J8").Value = 0.6429093297593833
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6351283333333334
$ws.Range("N8").Value = 1.905385
$ws.Range("O8").Value = 0.01806294655581008
$ws.Range("P8").Value = 0.01806294655581008
$ws.Range("Q8").Value = 3.931376212892222
$ws.Range("R8").Value = 35.38238591603
$ws.Range("S8").Value = 0.01161283686367542
$ws.Range("T8").Value = 0.01161283686367542

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.189892666666666
$ws.Range("H9").Value = 18.569678
$ws.Range("I9").Value = 0.6429093297593833
$ws.Range("J9").Value = 0.6429093297593833
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.516752666666667
$ws.Range("N9").Value = 28.550258
$ws.Range("O9").Value = 0.2706548988307292
$ws.Range("P9").Value = 0.2706548988307293
$ws.Range("Q9").Value = 58.90767754188045
$ws.Range("R9").Value = 530.169097876924
$ws.Range("S9").Value = 0.1740065596033578
$ws.Range("T9").Value = 0.1740065596033579

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.012114666666667
$ws.Range("H10").Value = 3.036344
$ws.Range("I10").Value = 0.105122656728831
$ws.Range("J10").Value = 0.105122656728831
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 24.244885
$ws.Range("N10").Value = 72.734655
$ws.Range("O10").Value = 0.6895205882382217
$ws.Range("P10").Value = 0.6895205882382218
$ws.Range("Q10").Value = 24.53860370014667
$ws.Range("R10").Value = 220.84743330132
$ws.Range("S10").Value = 0.07248423610482821
$ws.Range("T10").Value = 0.07248423610482822

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.012114666666667
$ws.Range("H11").Value = 3.036344
$ws.Range("I11").Value = 0.105122656728831
$ws.Range("J11").Value = 0.105122656728831
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7651789999999999
$ws.Range("N11").Value = 2.295537
$ws.Range("O11").Value = 0.02176156637523891
$ws.Range("P11").Value = 0.02176156637523891
$ws.Range("Q11").Value = 0.7744488885253333
$ws.Range("R11").Value = 6.970039996728
$ws.Range("S11").Value = 0.002287633671945911
$ws.Range("T11").Value = 0.002287633671945912

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.012114666666667
$ws.Range("H12").Value = 3.036344
$ws.Range("I12").Value = 0.105122656728831
$ws.Range("J12").Value = 0.105122656728831
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6351283333333334
$ws.Range("N12").Value = 1.905385
$ws.Range("O12").Value = 0.01806294655581008
$ws.Range("P12").Value = 0.01806294655581008
$ws.Range("Q12").Value = 0.6428227013822223
$ws.Range("R12").Value = 5.785404312440001
$ws.Range("S12").Value = 0.001898824930297643
$ws.Range("T12").Value = 0.001898824930297643

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.012114666666667
$ws.Range("H13").Value = 3.036344
$ws.Range("I13").Value = 0.105122656728831
$ws.Range("J13").Value = 0.105122656728831
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.516752666666667
$ws.Range("N13").Value = 28.550258
$ws.Range("O13").Value = 0.2706548988307292
$ws.Range("P13").Value = 0.2706548988307293
$ws.Range("Q13").Value = 9.632044952972445
$ws.Range("R13").Value = 86.68840457675201
$ws.Range("S13").Value = 0.02845196202175923
$ws.Range("T13").Value = 0.02845196202175924

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.008073333333333
$ws.Range("H14").Value = 3.02422
$ws.Range("I14").Value = 0.1047029061702051
$ws.Range("J14").Value = 0.1047029061702051
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 24.244885
$ws.Range("N14").Value = 72.734655
$ws.Range("O14").Value = 0.6895205882382217
$ws.Range("P14").Value = 0.6895205882382218
$ws.Range("Q14").Value = 24.44062203823334
$ws.Range("R14").Value = 219.9655983441
$ws.Range("S14").Value = 0.07219480945273117
$ws.Range("T14").Value = 0.07219480945273118

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.008073333333333
$ws.Range("H15").Value = 3.02422
$ws.Range("I15").Value = 0.1047029061702051
$ws.Range("J15").Value = 0.1047029061702051
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.7651789999999999
$ws.Range("N15").Value = 2.295537
$ws.Range("O15").Value = 0.02176156637523891
$ws.Range("P15").Value = 0.02176156637523891
$ws.Range("Q15").Value = 0.7713565451266666
$ws.Range("R15").Value = 6.94220890614
$ws.Range("S15").Value = 0.00227849924230333
$ws.Range("T15").Value = 0.002278499242303331

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.008073333333333
$ws.Range("H16").Value = 3.02422
$ws.Range("I16").Value = 0.1047029061702051
$ws.Range("J16").Value = 0.1047029061702051
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6351283333333334
$ws.Range("N16").Value = 1.905385
$ws.Range("O16").Value = 0.01806294655581008
$ws.Range("P16").Value = 0.01806294655581008
$ws.Range("Q16").Value = 0.6402559360777779
$ws.Range("R16").Value = 5.762303424700001
$ws.Range("S16").Value = 0.001891242998390412
$ws.Range("T16").Value = 0.001891242998390412

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.008073333333333
$ws.Range("H17").Value = 3.02422
$ws.Range("I17").Value = 0.1047029061702051
$ws.Range("J17").Value = 0.1047029061702051
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.516752666666667
$ws.Range("N17").Value = 28.550258
$ws.Range("O17").Value = 0.2706548988307292
$ws.Range("P17").Value = 0.2706548988307293
$ws.Range("Q17").Value = 9.593584583195556
$ws.Range("R17").Value = 86.34226124876
$ws.Range("S17").Value = 0.02833835447678021
$ws.Range("T17").Value = 0.02833835447678021

